# Refresh scraped schedule data for "horarios-141" workbook (new scrape at 08:41:16).
# For each sheet: update the "last updated" / "total rows" header cells, then apply
# the per-row cell changes coming from the new scrape merge (existing in-flight rows
# get refreshed Hora_Scrap/Hora_Llegada/Linea/Minutos, arrived rows are replaced, and
# new rows are appended at the bottom of the sorted-by-arrival-time table).
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: "LP1912" ----
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2,1).Value = "Última actualización: 08:41:16"
$ws.Cells.Item(3,1).Value = "Total filas: 91"
$ws.Cells.Item(37,1).Value = "07:28:23"
$ws.Cells.Item(37,3).Value = "16_SANTA ANA"
$ws.Cells.Item(37,4).Value = 3

$ws.Cells.Item(38,1).Value = "06:55:48"
$ws.Cells.Item(38,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(38,4).Value = 36

$ws.Cells.Item(60,1).Value = "08:41:16"
$ws.Cells.Item(60,2).Value = "08:41"
$ws.Cells.Item(60,4).Value = 0

$ws.Cells.Item(61,1).Value = "06:55:48"
$ws.Cells.Item(61,2).Value = "08:42"
$ws.Cells.Item(61,3).Value = "81_EL PELIGRO"
$ws.Cells.Item(61,4).Value = 107

$ws.Cells.Item(62,1).Value = "08:41:16"
$ws.Cells.Item(62,2).Value = "08:43"
$ws.Cells.Item(62,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(62,4).Value = 2

$ws.Cells.Item(63,2).Value = "08:43"
$ws.Cells.Item(63,3).Value = "14_ABASTO"
$ws.Cells.Item(63,4).Value = 39

$ws.Cells.Item(64,1).Value = "08:41:16"
$ws.Cells.Item(64,2).Value = "08:49"
$ws.Cells.Item(64,3).Value = "14_ABASTO"
$ws.Cells.Item(64,4).Value = 8

$ws.Cells.Item(65,1).Value = "07:28:23"
$ws.Cells.Item(65,2).Value = "08:50"
$ws.Cells.Item(65,3).Value = "81_EL PELIGRO"
$ws.Cells.Item(65,4).Value = 82

$ws.Cells.Item(66,1).Value = "08:41:16"
$ws.Cells.Item(66,2).Value = "08:53"
$ws.Cells.Item(66,3).Value = "10_OLMOS"
$ws.Cells.Item(66,4).Value = 12

$ws.Cells.Item(67,1).Value = "08:41:16"
$ws.Cells.Item(67,2).Value = "08:54"
$ws.Cells.Item(67,3).Value = "17_ROMERO"
$ws.Cells.Item(67,4).Value = 13

$ws.Cells.Item(68,1).Value = "08:41:16"
$ws.Cells.Item(68,2).Value = "09:01"
$ws.Cells.Item(68,3).Value = "215A_EL PATO"
$ws.Cells.Item(68,4).Value = 20

$ws.Cells.Item(69,1).Value = "08:04:39"
$ws.Cells.Item(69,2).Value = "09:02"
$ws.Cells.Item(69,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(69,4).Value = 58

$ws.Cells.Item(70,2).Value = "09:03"
$ws.Cells.Item(70,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(70,4).Value = 59

$ws.Cells.Item(71,1).Value = "08:41:16"
$ws.Cells.Item(71,2).Value = "09:04"
$ws.Cells.Item(71,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(71,4).Value = 23

$ws.Cells.Item(72,1).Value = "08:41:16"
$ws.Cells.Item(72,2).Value = "09:06"
$ws.Cells.Item(72,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(72,4).Value = 25

$ws.Cells.Item(73,2).Value = "09:10"
$ws.Cells.Item(73,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(73,4).Value = 66

$ws.Cells.Item(74,1).Value = "08:41:16"
$ws.Cells.Item(74,2).Value = "09:11"
$ws.Cells.Item(74,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(74,4).Value = 30

$ws.Cells.Item(75,1).Value = "08:41:16"
$ws.Cells.Item(75,2).Value = "09:11"
$ws.Cells.Item(75,3).Value = "16_SANTA ANA"
$ws.Cells.Item(75,4).Value = 30

$ws.Cells.Item(76,2).Value = "09:16"
$ws.Cells.Item(76,3).Value = "27_EL RETIRO"
$ws.Cells.Item(76,4).Value = 72

$ws.Cells.Item(77,1).Value = "08:41:16"
$ws.Cells.Item(77,2).Value = "09:17"
$ws.Cells.Item(77,3).Value = "27_EL RETIRO"
$ws.Cells.Item(77,4).Value = 36

$ws.Cells.Item(78,1).Value = "08:41:16"
$ws.Cells.Item(78,2).Value = "09:21"
$ws.Cells.Item(78,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(78,4).Value = 40

$ws.Cells.Item(79,1).Value = "08:04:39"
$ws.Cells.Item(79,2).Value = "09:22"
$ws.Cells.Item(79,3).Value = "17_ROMERO"
$ws.Cells.Item(79,4).Value = 78
$ws.Cells.Item(79,5).Value = "LP1912"

$ws.Cells.Item(80,1).Value = "08:04:39"
$ws.Cells.Item(80,2).Value = "09:23"
$ws.Cells.Item(80,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(80,4).Value = 79
$ws.Cells.Item(80,5).Value = "LP1912"

$ws.Cells.Item(81,1).Value = "08:41:16"
$ws.Cells.Item(81,2).Value = "09:23"
$ws.Cells.Item(81,3).Value = "17_ROMERO"
$ws.Cells.Item(81,4).Value = 42
$ws.Cells.Item(81,5).Value = "LP1912"

$ws.Cells.Item(82,1).Value = "08:41:16"
$ws.Cells.Item(82,2).Value = "09:23"
$ws.Cells.Item(82,3).Value = "16_SANTA ANA"
$ws.Cells.Item(82,4).Value = 42
$ws.Cells.Item(82,5).Value = "LP1912"

$ws.Cells.Item(83,1).Value = "08:41:16"
$ws.Cells.Item(83,2).Value = "09:24"
$ws.Cells.Item(83,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(83,4).Value = 43
$ws.Cells.Item(83,5).Value = "LP1912"

$ws.Cells.Item(84,1).Value = "08:04:39"
$ws.Cells.Item(84,2).Value = "09:25"
$ws.Cells.Item(84,3).Value = "81_EL PELIGRO"
$ws.Cells.Item(84,4).Value = 81
$ws.Cells.Item(84,5).Value = "LP1912"

$ws.Cells.Item(85,1).Value = "08:41:16"
$ws.Cells.Item(85,2).Value = "09:32"
$ws.Cells.Item(85,3).Value = "15_ABASTO"
$ws.Cells.Item(85,4).Value = 51
$ws.Cells.Item(85,5).Value = "LP1912"

$ws.Cells.Item(86,1).Value = "08:41:16"
$ws.Cells.Item(86,2).Value = "09:33"
$ws.Cells.Item(86,3).Value = "10_OLMOS"
$ws.Cells.Item(86,4).Value = 52
$ws.Cells.Item(86,5).Value = "LP1912"

$ws.Cells.Item(87,1).Value = "08:41:16"
$ws.Cells.Item(87,2).Value = "09:35"
$ws.Cells.Item(87,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(87,4).Value = 54
$ws.Cells.Item(87,5).Value = "LP1912"

$ws.Cells.Item(88,1).Value = "08:04:39"
$ws.Cells.Item(88,2).Value = "09:41"
$ws.Cells.Item(88,3).Value = "215C_EL PATO"
$ws.Cells.Item(88,4).Value = 97
$ws.Cells.Item(88,5).Value = "LP1912"

$ws.Cells.Item(89,1).Value = "08:41:16"
$ws.Cells.Item(89,2).Value = "09:42"
$ws.Cells.Item(89,3).Value = "215C_EL PATO"
$ws.Cells.Item(89,4).Value = 61
$ws.Cells.Item(89,5).Value = "LP1912"

$ws.Cells.Item(90,1).Value = "08:04:39"
$ws.Cells.Item(90,2).Value = "09:43"
$ws.Cells.Item(90,3).Value = "14_ABASTO"
$ws.Cells.Item(90,4).Value = 99
$ws.Cells.Item(90,5).Value = "LP1912"

$ws.Cells.Item(91,1).Value = "08:41:16"
$ws.Cells.Item(91,2).Value = "09:44"
$ws.Cells.Item(91,3).Value = "14_ABASTO"
$ws.Cells.Item(91,4).Value = 63
$ws.Cells.Item(91,5).Value = "LP1912"

$ws.Cells.Item(92,1).Value = "08:41:16"
$ws.Cells.Item(92,2).Value = "09:52"
$ws.Cells.Item(92,3).Value = "15_ABASTO"
$ws.Cells.Item(92,4).Value = 71
$ws.Cells.Item(92,5).Value = "LP1912"

$ws.Cells.Item(93,1).Value = "08:41:16"
$ws.Cells.Item(93,2).Value = "09:53"
$ws.Cells.Item(93,3).Value = "10_OLMOS"
$ws.Cells.Item(93,4).Value = 72
$ws.Cells.Item(93,5).Value = "LP1912"

$ws.Cells.Item(94,1).Value = "08:41:16"
$ws.Cells.Item(94,2).Value = "10:11"
$ws.Cells.Item(94,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(94,4).Value = 90
$ws.Cells.Item(94,5).Value = "LP1912"

$ws.Cells.Item(95,1).Value = "08:41:16"
$ws.Cells.Item(95,2).Value = "10:21"
$ws.Cells.Item(95,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(95,4).Value = 100
$ws.Cells.Item(95,5).Value = "LP1912"

$ws.Cells.Item(96,1).Value = "08:41:16"
$ws.Cells.Item(96,2).Value = "10:26"
$ws.Cells.Item(96,3).Value = "215A_EL PATO"
$ws.Cells.Item(96,4).Value = 105
$ws.Cells.Item(96,5).Value = "LP1912"


# ---- Sheet 2: "LP1912-215" ----
$ws = $wb.Worksheets.Item(2)

$ws.Cells.Item(2,1).Value = "Última actualización: 08:41:16"
$ws.Cells.Item(3,1).Value = "Total filas: 12"
$ws.Cells.Item(14,1).Value = "08:41:16"
$ws.Cells.Item(14,4).Value = 20

$ws.Cells.Item(16,1).Value = "08:41:16"
$ws.Cells.Item(16,2).Value = "09:42"
$ws.Cells.Item(16,3).Value = "215C_EL PATO"
$ws.Cells.Item(16,4).Value = 61
$ws.Cells.Item(16,5).Value = "LP1912"

$ws.Cells.Item(17,1).Value = "08:41:16"
$ws.Cells.Item(17,2).Value = "10:26"
$ws.Cells.Item(17,3).Value = "215A_EL PATO"
$ws.Cells.Item(17,4).Value = 105
$ws.Cells.Item(17,5).Value = "LP1912"


# ---- Sheet 3: "6203-6173" ----
$ws = $wb.Worksheets.Item(3)

$ws.Cells.Item(2,1).Value = "Última actualización: 08:41:16"
$ws.Cells.Item(3,1).Value = "Total filas: 19"
$ws.Cells.Item(22,1).Value = "08:41:16"
$ws.Cells.Item(22,2).Value = "09:10"
$ws.Cells.Item(22,3).Value = "215D_LA PLATA"
$ws.Cells.Item(22,4).Value = 29
$ws.Cells.Item(22,5).Value = "L6203"

$ws.Cells.Item(23,1).Value = "08:04:39"
$ws.Cells.Item(23,2).Value = "10:02"
$ws.Cells.Item(23,3).Value = "215B_LP-P MOR-40 Y 115"
$ws.Cells.Item(23,4).Value = 118
$ws.Cells.Item(23,5).Value = "L6173"

$ws.Cells.Item(24,1).Value = "08:41:16"
$ws.Cells.Item(24,2).Value = "10:03"
$ws.Cells.Item(24,3).Value = "215B_LP-P MOR-40 Y 115"
$ws.Cells.Item(24,4).Value = 82
$ws.Cells.Item(24,5).Value = "L6173"

